# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Leve profit sheets
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H:N) across all 8 sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2158.5454
$ws.Range("I62").Value = 1749.3334
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 1749.3334
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -1125.3334
$ws.Range("N62").Value = -5248
# Row 65
$ws.Range("H65").Value = 2158.5454
$ws.Range("I65").Value = 1749.3334
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 8746.666999999999
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -5626.666999999999
$ws.Range("N65").Value = -26240
# Row 70
$ws.Range("H70").Value = 898
$ws.Range("I70").Value = 956.4
$ws.Range("J70").Value = 868.8
$ws.Range("K70").Value = 2869.2
$ws.Range("L70").Value = 2606.4
$ws.Range("M70").Value = -2599.2
$ws.Range("N70").Value = -3146.4
# Row 73
$ws.Range("H73").Value = 898
$ws.Range("I73").Value = 956.4
$ws.Range("J73").Value = 868.8
$ws.Range("K73").Value = 2869.2
$ws.Range("L73").Value = 2606.4
$ws.Range("M73").Value = -1933.2
$ws.Range("N73").Value = -4478.4
# Row 103
$ws.Range("H103").Value = 167191.5
$ws.Range("I103").Value = 250387.25
$ws.Range("J103").Value = 800
$ws.Range("K103").Value = 751161.75
$ws.Range("L103").Value = 2400
$ws.Range("M103").Value = -750575.75
$ws.Range("N103").Value = -3572
# Row 116
$ws.Range("H116").Value = 2614.5454
$ws.Range("I116").Value = 2488.75
$ws.Range("K116").Value = 2488.75
$ws.Range("M116").Value = 953.25
# Row 137
$ws.Range("H137").Value = 4029.9333
$ws.Range("I137").Value = 4411.6665
$ws.Range("J137").Value = 2503
$ws.Range("K137").Value = 13234.9995
$ws.Range("L137").Value = 7509
$ws.Range("M137").Value = -10684.9995
$ws.Range("N137").Value = -12609
# Row 138
$ws.Range("H138").Value = 238872.48
$ws.Range("I138").Value = 2296.3684
$ws.Range("J138").Value = 405351.97
$ws.Range("K138").Value = 6889.1052
$ws.Range("L138").Value = 1216055.91
$ws.Range("M138").Value = -1749.1052
$ws.Range("N138").Value = -1226335.91

$ws = $wb.Worksheets.Item("ARM")
# Row 15
$ws.Range("H15").Value = 69000
$ws.Range("J15").Value = 69000
$ws.Range("L15").Value = 69000
$ws.Range("N15").Value = -69700
# Row 32
$ws.Range("H32").Value = 1016163.8
$ws.Range("I32").Value = 1206406.2
$ws.Range("J32").Value = 26902.8
$ws.Range("K32").Value = 1206406.2
$ws.Range("L32").Value = 26902.8
$ws.Range("M32").Value = -1206119.2
$ws.Range("N32").Value = -27476.8
# Row 61
$ws.Range("H61").Value = 2436.9048
$ws.Range("I61").Value = 1685.0625
$ws.Range("K61").Value = 1685.0625
$ws.Range("M61").Value = -1473.0625
# Row 88
$ws.Range("H88").Value = 2475.6843
$ws.Range("I88").Value = 2474.3333
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 2474.3333
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -2068.3333
$ws.Range("N88").Value = -3312
# Row 91
$ws.Range("H91").Value = 2475.6843
$ws.Range("I91").Value = 2474.3333
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 2474.3333
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -1070.3333
$ws.Range("N91").Value = -5308
# Row 136
$ws.Range("H136").Value = 2436.9048
$ws.Range("I136").Value = 1685.0625
$ws.Range("K136").Value = 5055.1875
$ws.Range("M136").Value = -2505.1875

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 90911610
$ws.Range("I86").Value = 125001840
$ws.Range("K86").Value = 125001840
$ws.Range("M86").Value = -125000717
# Row 89
$ws.Range("H89").Value = 90911610
$ws.Range("I89").Value = 125001840
$ws.Range("K89").Value = 625009200
$ws.Range("M89").Value = -625003584
# Row 94
$ws.Range("H94").Value = 602.5625
$ws.Range("I94").Value = 525.2
$ws.Range("K94").Value = 525.2
$ws.Range("M94").Value = -74.20000000000005

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3598.6924
$ws.Range("I31").Value = 1074.5
$ws.Range("J31").Value = 7637.4
$ws.Range("K31").Value = 1074.5
$ws.Range("L31").Value = 7637.4
$ws.Range("M31").Value = -779.5
$ws.Range("N31").Value = -8227.4
# Row 34
$ws.Range("H34").Value = 3598.6924
$ws.Range("I34").Value = 1074.5
$ws.Range("J34").Value = 7637.4
$ws.Range("K34").Value = 1074.5
$ws.Range("L34").Value = 7637.4
$ws.Range("M34").Value = -872.5
$ws.Range("N34").Value = -8041.4
# Row 62
$ws.Range("H62").Value = 4101.643
$ws.Range("I62").Value = 4102.3
$ws.Range("J62").Value = 4100
$ws.Range("K62").Value = 4102.3
$ws.Range("L62").Value = 4100
$ws.Range("M62").Value = -3478.3
$ws.Range("N62").Value = -5348
# Row 65
$ws.Range("H65").Value = 4101.643
$ws.Range("I65").Value = 4102.3
$ws.Range("J65").Value = 4100
$ws.Range("K65").Value = 20511.5
$ws.Range("L65").Value = 20500
$ws.Range("M65").Value = -17391.5
$ws.Range("N65").Value = -26740
# Row 122
$ws.Range("H122").Value = 1493
$ws.Range("I122").Value = 1079.6875
$ws.Range("J122").Value = 1747.3462
$ws.Range("K122").Value = 3239.0625
$ws.Range("L122").Value = 5242.0386
$ws.Range("M122").Value = -789.0625
$ws.Range("N122").Value = -10142.0386
# Row 132
$ws.Range("H132").Value = 6668521.5
$ws.Range("I132").Value = 1285.1875
$ws.Range("J132").Value = 18521386
$ws.Range("K132").Value = 3855.5625
$ws.Range("L132").Value = 55564158
$ws.Range("M132").Value = -1325.5625
$ws.Range("N132").Value = -55569218

$ws = $wb.Worksheets.Item("CUL")
# Row 58
$ws.Range("H58").Value = 1735.6522
$ws.Range("J58").Value = 2058.8235
$ws.Range("L58").Value = 6176.470499999999
$ws.Range("N58").Value = -6432.470499999999
# Row 122
$ws.Range("H122").Value = 7471.533
$ws.Range("I122").Value = 567.5
$ws.Range("J122").Value = 21279.6
$ws.Range("K122").Value = 5107.5
$ws.Range("L122").Value = 191516.4
$ws.Range("M122").Value = -2657.5
$ws.Range("N122").Value = -196416.4
# Row 131
$ws.Range("H131").Value = 943.069
$ws.Range("I131").Value = 290
$ws.Range("J131").Value = 1047.56
$ws.Range("K131").Value = 870
$ws.Range("L131").Value = 3142.68
$ws.Range("M131").Value = 4170
$ws.Range("N131").Value = -13222.68
# Row 137
$ws.Range("H137").Value = 9295.467000000001
$ws.Range("J137").Value = 3258.25
$ws.Range("L137").Value = 9774.75
$ws.Range("N137").Value = -19974.75
# Row 140
$ws.Range("H140").Value = 2376.5789
$ws.Range("I140").Value = 1493.75
$ws.Range("J140").Value = 3890
$ws.Range("K140").Value = 4481.25
$ws.Range("L140").Value = 11670
$ws.Range("M140").Value = 698.75
$ws.Range("N140").Value = -22030

$ws = $wb.Worksheets.Item("GSM")
# Row 12
$ws.Range("H12").Value = 70004
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 70004
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 70004
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -70284
# Row 70
$ws.Range("H70").Value = 5249.8394
$ws.Range("I70").Value = 5087.486
$ws.Range("J70").Value = 5520.4287
$ws.Range("K70").Value = 5087.486
$ws.Range("L70").Value = 5520.4287
$ws.Range("M70").Value = -4817.486
$ws.Range("N70").Value = -6060.4287
# Row 73
$ws.Range("H73").Value = 5249.8394
$ws.Range("I73").Value = 5087.486
$ws.Range("J73").Value = 5520.4287
$ws.Range("K73").Value = 5087.486
$ws.Range("L73").Value = 5520.4287
$ws.Range("M73").Value = -4151.486
$ws.Range("N73").Value = -7392.4287
# Row 132
$ws.Range("H132").Value = 1824.84
$ws.Range("I132").Value = 1146.7646
$ws.Range("J132").Value = 3265.75
$ws.Range("K132").Value = 3440.2938
$ws.Range("L132").Value = 9797.25
$ws.Range("M132").Value = -910.2937999999999
$ws.Range("N132").Value = -14857.25

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 59502
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 59502
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 59502
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -59726
# Row 7
$ws.Range("H7").Value = 66669760
$ws.Range("I7").Value = 83336630
$ws.Range("J7").Value = 2268.3333
$ws.Range("K7").Value = 83336630
$ws.Range("L7").Value = 2268.3333
$ws.Range("M7").Value = -83336518
$ws.Range("N7").Value = -2492.3333
# Row 126
$ws.Range("H126").Value = 66669760
$ws.Range("I126").Value = 83336630
$ws.Range("J126").Value = 2268.3333
$ws.Range("K126").Value = 250009890
$ws.Range("L126").Value = 6804.999899999999
$ws.Range("M126").Value = -250007420
$ws.Range("N126").Value = -11744.9999
# Row 132
$ws.Range("H132").Value = 3451.2166
$ws.Range("I132").Value = 3161.2778
$ws.Range("K132").Value = 9483.8334
$ws.Range("M132").Value = -6953.8334
# Row 139
$ws.Range("H139").Value = 9695000
$ws.Range("J139").Value = 90000
$ws.Range("L139").Value = 90000
$ws.Range("N139").Value = -100280

$ws = $wb.Worksheets.Item("WVR")
# Row 82
$ws.Range("H82").Value = 39980
$ws.Range("J82").Value = 39980
$ws.Range("L82").Value = 39980
$ws.Range("N82").Value = -40746
# Row 85
$ws.Range("H85").Value = 39980
$ws.Range("J85").Value = 39980
$ws.Range("L85").Value = 39980
$ws.Range("N85").Value = -42632
# Row 122
$ws.Range("H122").Value = 2303.0527
$ws.Range("I122").Value = 2198
$ws.Range("K122").Value = 6594
$ws.Range("M122").Value = -4144
# Row 126
$ws.Range("H126").Value = 1354.1666
$ws.Range("I126").Value = 1081.25
$ws.Range("K126").Value = 3243.75
$ws.Range("M126").Value = -773.75
# Row 132
$ws.Range("H132").Value = 7940284
$ws.Range("I132").Value = 5722.625
$ws.Range("J132").Value = 12823091
$ws.Range("K132").Value = 17167.875
$ws.Range("L132").Value = 38469273
$ws.Range("M132").Value = -14637.875
$ws.Range("N132").Value = -38474333
# Row 139
$ws.Range("H139").Value = 94378.42999999999
$ws.Range("J139").Value = 94378.42999999999
$ws.Range("L139").Value = 94378.42999999999
$ws.Range("N139").Value = -104658.43
